$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tbl = $sh.Table

# Change the table style to the new built-in style.
$tbl.ApplyStyle("{7CF10278-3543-4173-98FF-44D05C8B350D}")

# Update the Twitter followers count from 142 to 144.
$cell = $tbl.Cell(1, 2)
$cell.Shape.TextFrame.TextRange.Text = "144"
